# "Refined metadata to be additional tab"
#
# 1. Refresh the panel-query timestamps (time_taken, column F) on the
#    existing "data" sheet.
# 2. Add a new "metadata" worksheet after "data" holding the panel-level
#    metadata (name/id/version/etc) that used to live only implicitly.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Update the time_taken timestamps on the "data" sheet ---------------
$dataSheet.Range("F2").Value = "2021-10-05 14:21:36.994118"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:36.994126"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:36.994130"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:36.994133"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:36.994137"

# --- 2. Add the new "metadata" worksheet, placed right after "data" --------
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Header row (B1:G1)
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Monogenic nephrogenic diabetes insipidus"
$newSheet.Range("C2").Value = 18

# data_version ("1.8") must stay a text value, not be coerced to a number.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.8"
$newSheet.Range("D2").ClearFormats()

$newSheet.Range("E2").Value = "2018-11-19T12:28:13.664834Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:36.990533"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/18/?format=json"

# Match the header styling (bold / centered / bordered) already used by the
# "data" sheet's own header row, and the index-column style used for A2.
$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$dataSheet.Activate()
$dataSheet.Range("A1").Select()
